# Insert a new data row after row 112 (i.e. before the former row 113),
# shifting all subsequent rows down by one, then populate the new row
# with the added price-record for Jengibre (Vega Modelo de Temuco).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("113:113").Insert()

$ws.Range("A113").Value = 10
$ws.Range("B113").Value = 'Vega Modelo de Temuco'
$ws.Range("C113").Value = 'La Araucanía'
$ws.Range("D113").Value = 44907
$ws.Range("E113").Value = 9
$ws.Range("F113").Value = 100114007
$ws.Range("G113").Value = 'Jengibre'
$ws.Range("H113").Value = 'Sin especificar'
$ws.Range("I113").Value = 'Primera'
$ws.Range("J113").Value = 70
$ws.Range("K113").Value = 20000
$ws.Range("L113").Value = 20000
$ws.Range("M113").Value = 20000
$ws.Range("N113").Value = '$/caja 13 kilos'
$ws.Range("O113").Value = 'Perú'
$ws.Range("P113").Value = 1538
$ws.Range("Q113").Value = 13
$ws.Range("R113").Value = 'Hortaliza'
